$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G: numeric header (row 1), text header "access" (row 2),
# and a value of 1 for every data row (3..251) — mirrors the existing
# A:F column layout (numbered header row, label header row, data rows).

# Row 1 header number — match formatting of the existing numeric header cells (A1:F1)
$ws.Range("G1").Value = 7
$ws.Range("A1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null

# Row 2 column title — match formatting of the existing title cells (A2:F2)
$ws.Range("G2").Value = "access"
$ws.Range("A2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Data rows: flag every record with a 1 in the new "access" column
for ($r = 3; $r -le 251; $r++) {
    $ws.Cells.Item($r, 7).Value = 1
}

$ws.Range("G8").Select() | Out-Null
